$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab from "Evaluation" to "Sheet1"
$ws.Name = "Sheet1"

# Row 15: replace the Cohere translation text and its BLEU score
$ws.Range("C15").Value = "ஒரு முக்கிய கிறித்தவ நம்பிக்கை என்பது, செல்வாக்கு துன்பத்தையும் வறுமையையும் நீக்குவதற்காக பயன்படுத்தப்பட வேண்டும் என்பதும், அதற்காகவே தேவாலயத்தின் நிதி இருப்பதாகும் என்பது ஆகும்."
$ws.Range("D15").Value = 0.011401178105196599

# Row 80: replace the Cohere translation text and its BLEU score
$ws.Range("C80").Value = "இவசகியின் பயணத்தின்போது, அவர் பல சந்தர்ப்பங்களில் சிக்கலில் மாட்டிக் கொண்டார்."
$ws.Range("D80").Value = 0.078825884423426515
